# disable-showcase-macro.xlsx — add the JSON `storeKeys(json,jsonpath,var)` macro
# ("keys(jsonpath)" expression) to the hidden '#system' catalogue sheet and keep
# the workbook-level defined names (used by the showcase UI) in sync.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------------
# 1) JSON category (column M): the per-category value lists are kept in
#    alphabetical order. "storeKeys(...)" sorts between "minify(...)" (M15)
#    and "storeValue(...)" (M16), so push M16:M17 down one row and write the
#    new entry into the freed M16 slot. This grows the `json` named range
#    from M2:M17 to M2:M18.
# ---------------------------------------------------------------------------
$ws.Range("M16").Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)
$ws.Range("M16").Value = "storeKeys(json,jsonpath,var)"

# ---------------------------------------------------------------------------
# 2) target list (column A): column A enumerates every category name
#    (alphabetically) for the picker UI. The "text" category (A25) is no
#    longer exposed there, so remove that single cell and shift A26:A31 up.
#    This shrinks the `target` named range from A2:A31 to A2:A30.
# ---------------------------------------------------------------------------
$ws.Range("A25").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

# ---------------------------------------------------------------------------
# 3) "text" column (Y): its data column is removed outright, so every
#    category to its right (web, webalert, webcookie, ws, ws.async, xml)
#    shifts one column to the left (Z->Y, AA->Z, AB->AA, AC->AB, AD->AC,
#    AE->AD).
# ---------------------------------------------------------------------------
$ws.Columns.Item("Y").Delete()

# ---------------------------------------------------------------------------
# 4) Keep the workbook-level defined names lined up with the moved ranges.
# ---------------------------------------------------------------------------
$wb.Names.Item("json").RefersTo = "='#system'!`$M`$2:`$M`$18"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("web").RefersTo = "='#system'!`$Y`$2:`$Y`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AD`$2:`$AD`$27"
